# Apply the "Заказ №1" packing list update:
#  - row 2 (article E3-55-Al-1000-4-pt3.0) becomes E3-55-Al-5000-3-uv, an
#    angled vertical section sized 450х450 instead of a straight 3000 section,
#    with an updated mass;
#  - row 3 (the second straight-section line) is removed entirely;
#  - the page header date is bumped from 11 10 2021 to 21 10 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 in place ------------------------------------------------
$ws.Range("B2").Value = "E3-55-Al-5000-3-uv"
$ws.Range("D2").Value = "угловая вертикальная секция"
$ws.Range("E2").Value = "450х450"
$ws.Range("I2").Value = 41.94

# --- Remove row 3 (second line item) --------------------------------------
$ws.Rows(3).Delete()

# --- Bump the date printed in the page header's right section ------------
$ws.PageSetup.RightHeader = $ws.PageSetup.RightHeader -replace "11 10 2021", "21 10 2021"
